# Apply updates to column F (想去人数) values across sheets, per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 34
$ws.Cells.Item(3, 6).Value = 8863
$ws.Cells.Item(4, 6).Value = 1955
$ws.Cells.Item(5, 6).Value = 6568
$ws.Cells.Item(6, 6).Value = 168
$ws.Cells.Item(7, 6).Value = 2118
$ws.Cells.Item(9, 6).Value = 67
$ws.Cells.Item(16, 6).Value = 8760
$ws.Cells.Item(28, 6).Value = 198
$ws.Cells.Item(31, 6).Value = 40
$ws.Cells.Item(34, 6).Value = 28
$ws.Cells.Item(35, 6).Value = 2221
$ws.Cells.Item(36, 6).Value = 870
$ws.Cells.Item(41, 6).Value = 244
$ws.Cells.Item(42, 6).Value = 176
$ws.Cells.Item(44, 6).Value = 817
$ws.Cells.Item(45, 6).Value = 82
$ws.Cells.Item(48, 6).Value = 3990

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 3
$ws.Cells.Item(12, 6).Value = 12

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 2336
$ws.Cells.Item(3, 6).Value = 716
$ws.Cells.Item(4, 6).Value = 323

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 2336
$ws.Cells.Item(3, 6).Value = 716
$ws.Cells.Item(4, 6).Value = 34
$ws.Cells.Item(5, 6).Value = 8863
$ws.Cells.Item(7, 6).Value = 323
$ws.Cells.Item(8, 6).Value = 1955
$ws.Cells.Item(9, 6).Value = 2118
$ws.Cells.Item(11, 6).Value = 67
$ws.Cells.Item(16, 6).Value = 12
$ws.Cells.Item(19, 6).Value = 8760
$ws.Cells.Item(28, 6).Value = 198
$ws.Cells.Item(33, 6).Value = 28
$ws.Cells.Item(34, 6).Value = 2221
$ws.Cells.Item(35, 6).Value = 870
$ws.Cells.Item(39, 6).Value = 244
$ws.Cells.Item(41, 6).Value = 176
$ws.Cells.Item(44, 6).Value = 3990
